# Edit script: implements the selectable-editing change across all lists.
# Updates existing rows 2-8 values and appends new rows 9-14 to the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all data cells (existing + new rows) keep their original text
# representation (e.g. "4000.00") instead of being auto-converted to
# numbers by Excel when the value looks numeric.
$ws.Range("A2:G14").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "31/03/2001"
$ws.Range("B2").Value = "4000.00"
$ws.Range("C2").Value = "4000.00"
$ws.Range("D2").Value = "4000.00"
$ws.Range("E2").Value = "4000.00"
$ws.Range("F2").Value = "0.00"
$ws.Range("G2").Value = "100.00"

# Row 3 (A, B, D, F, G unchanged)
$ws.Range("C3").Value = "5000.00"
$ws.Range("E3").Value = "5000.00"

# Row 4
$ws.Range("A4").Value = "31/04/2000"
$ws.Range("B4").Value = "4000.00"
$ws.Range("C4").Value = "9000.00"
$ws.Range("D4").Value = "4000.00"
$ws.Range("E4").Value = "9000.00"
$ws.Range("F4").Value = "0.00"
$ws.Range("G4").Value = "100.00"

# Row 5 (A, B, D, F, G unchanged)
$ws.Range("C5").Value = "10000.00"
$ws.Range("E5").Value = "10000.00"

# Row 6 (A, B, D, F, G unchanged)
$ws.Range("C6").Value = "12000.00"
$ws.Range("E6").Value = "12000.00"

# Row 7 (A, B, D, F, G unchanged)
$ws.Range("C7").Value = "14000.00"
$ws.Range("E7").Value = "14000.00"

# Row 8 (A, F unchanged)
$ws.Range("B8").Value = "6000.00"
$ws.Range("C8").Value = "20000.00"
$ws.Range("D8").Value = "6000.00"
$ws.Range("E8").Value = "20000.00"
$ws.Range("G8").Value = "100.00"

# New row 9
$ws.Range("A9").Value = "31/03/2001"
$ws.Range("B9").Value = "6000.00"
$ws.Range("C9").Value = "26000.00"
$ws.Range("D9").Value = "6000.00"
$ws.Range("E9").Value = "26000.00"
$ws.Range("F9").Value = "0.00"
$ws.Range("G9").Value = "100.00"

# New row 10
$ws.Range("A10").Value = "31/03/2001"
$ws.Range("B10").Value = "4000.00"
$ws.Range("C10").Value = "30000.00"
$ws.Range("D10").Value = "4000.00"
$ws.Range("E10").Value = "30000.00"
$ws.Range("F10").Value = "0.00"
$ws.Range("G10").Value = "100.00"

# New row 11
$ws.Range("A11").Value = "31/03/2001"
$ws.Range("B11").Value = "4000.00"
$ws.Range("C11").Value = "34000.00"
$ws.Range("D11").Value = "4000.00"
$ws.Range("E11").Value = "34000.00"
$ws.Range("F11").Value = "0.00"
$ws.Range("G11").Value = "100.00"

# New row 12
$ws.Range("A12").Value = "31/03/2001"
$ws.Range("B12").Value = "4000.00"
$ws.Range("C12").Value = "38000.00"
$ws.Range("D12").Value = "4000.00"
$ws.Range("E12").Value = "38000.00"
$ws.Range("F12").Value = "0.00"
$ws.Range("G12").Value = "100.00"

# New row 13
$ws.Range("A13").Value = "31/03/2001"
$ws.Range("B13").Value = "4000.00"
$ws.Range("C13").Value = "42000.00"
$ws.Range("D13").Value = "4000.00"
$ws.Range("E13").Value = "42000.00"
$ws.Range("F13").Value = "0.00"
$ws.Range("G13").Value = "100.00"

# New row 14
$ws.Range("A14").Value = "31/03/2001"
$ws.Range("B14").Value = "4000.00"
$ws.Range("C14").Value = "46000.00"
$ws.Range("D14").Value = "4000.00"
$ws.Range("E14").Value = "46000.00"
$ws.Range("F14").Value = "0.00"
$ws.Range("G14").Value = "100.00"
